$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: bold "abstractions" and "productivity" in the first sentence.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("abstractions") | Out-Null
$rng.Bold = 1

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("productivity") | Out-Null
$rng.Bold = 1

# ---------------------------------------------------------------------------
# Change 2: bold "art" and "science" in "computer architecture is the art
# and science of tradeoffs".
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("art") | Out-Null
$rng.Bold = 1

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("science") | Out-Null
$rng.Bold = 1

# ---------------------------------------------------------------------------
# Change 3: bold "Even the Von Neumann model builds an abstraction layer
# between the processors and memories" (without the trailing period) and
# move the _GoBack bookmark (currently at the end of the document) to wrap
# around this sentence.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("Even the Von Neumann model builds an abstraction layer between the processors and memories") | Out-Null
$rng.Bold = 1

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("Even the Von Neumann model builds an abstraction layer between the processors and memories") | Out-Null
$d.Bookmarks.Add("_GoBack", $rng) | Out-Null

# ---------------------------------------------------------------------------
# Change 4: remove the QoS spell-check proof errors by collapsing the whole
# sentence back into plain, un-annotated text (the stray bookmark that used
# to sit here was already relocated above).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute(" Also, energy consumption for refresh is a downside. Additionally refreshing DRAM periodically cause performance degradation (rows get unavailable while refreshing them), and QoS/ predictability problem because long pauses for refreshes decreases the QoS.", $true, $false, $false, $false, $false, $true, 1, $false, " Also, energy consumption for refresh is a downside. Additionally refreshing DRAM periodically cause performance degradation (rows get unavailable while refreshing them), and QoS/ predictability problem because long pauses for refreshes decreases the QoS.", 2) | Out-Null
